$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $inner = $val -replace '"', '""'
    $ws.Range($addr).Formula = ('="' + $inner + '"')
    $ws.Range($addr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null
}

Set-TextValue "D2" "27.506.34"
Set-TextValue "E2" "  -0.75%  "
Set-TextValue "D3" "1.831.17"
Set-TextValue "E3" "  -0.86%  "
Set-TextValue "E4" "  +0.03%  "
Set-TextValue "D5" "312.49"
Set-TextValue "D6" "1.002"
Set-TextValue "E6" "  +0.04%  "
Set-TextValue "D7" "0.4291"
Set-TextValue "E7" "  -0.63%  "
Set-TextValue "D8" "0.3660"
Set-TextValue "E8" "  +0.14%  "
Set-TextValue "E9" "  -0.76%  "
Set-TextValue "D10" "0.8639"
Set-TextValue "E10" "  -1.70%  "
Set-TextValue "D12" "1.948.92"
Set-TextValue "E12" "  +6.47%  "
Set-TextValue "D13" "5.395"
Set-TextValue "E13" "  +0.72%  "
Set-TextValue "D14" "6.535"
Set-TextValue "E14" "  -0.08%  "
Set-TextValue "D15" "0.06935"
Set-TextValue "E15" "  -0.47%  "
Set-TextValue "E16" "  -0.04%  "
Set-TextValue "D17" "80.74"
Set-TextValue "E17" "  +1.07%  "
Set-TextValue "D18" "0.000008902"
Set-TextValue "E18" "  -1.20%  "
Set-TextValue "D19" "1.001"
Set-TextValue "E19" "  +0.05%  "
Set-TextValue "D20" "15.39"
Set-TextValue "E20" "  -0.02%  "
Set-TextValue "D21" "27.715.41"
Set-TextValue "E21" "  +0.20%  "
Set-TextValue "D22" "5.145"
Set-TextValue "E22" "  +3.32%  "
Set-TextValue "D23" "10.84"
Set-TextValue "E23" "  +4.99%  "
Set-TextValue "D24" "2.113.90"
Set-TextValue "E24" "  +0.74%  "
Set-TextValue "D25" "1.994"
Set-TextValue "E25" "  +0.14%  "
Set-TextValue "D26" "154.43"
Set-TextValue "D27" "18.88"
Set-TextValue "E27" "  +1.30%  "
Set-TextValue "D28" "5.103"
Set-TextValue "E28" "  -2.91%  "
Set-TextValue "D29" "114.26"
Set-TextValue "E29" "  -4.77%  "
Set-TextValue "D30" "1.828"
Set-TextValue "E30" "  -3.12%  "
Set-TextValue "D31" "0.08850"
Set-TextValue "E31" "  -0.46%  "
Set-TextValue "D32" "0.7482"
Set-TextValue "E32" "  -1.50%  "
Set-TextValue "D33" "2.988"
Set-TextValue "E33" "  +0.80%  "
Set-TextValue "D34" "4.541"
Set-TextValue "E34" "  -0.29%  "
Set-TextValue "D35" "1.129"
Set-TextValue "E35" "  +0.32%  "
Set-TextValue "E36" "  +0.04%  "
Set-TextValue "D37" "1.087"
Set-TextValue "E37" "  -2.02%  "
Set-TextValue "D38" "0.05327"
Set-TextValue "E38" "  -2.15%  "
Set-TextValue "D39" "0.01934"
Set-TextValue "E39" "  -0.14%  "
Set-TextValue "D40" "2.796"
Set-TextValue "E40" "  -1.65%  "
Set-TextValue "D41" "0.5075"
Set-TextValue "E41" "  -0.42%  "
Set-TextValue "D42" "0.1660"
Set-TextValue "E42" "  -0.33%  "
Set-TextValue "D43" "6.525"
Set-TextValue "E43" "  -1.41%  "
Set-TextValue "D44" "8.309"
Set-TextValue "E44" "  -1.02%  "
Set-TextValue "D45" "10.43"
Set-TextValue "E45" "  +0.69%  "
Set-TextValue "D46" "105.88"
Set-TextValue "E46" "  +0.29%  "
Set-TextValue "D47" "0.06496"
Set-TextValue "E47" "  -0.93%  "
Set-TextValue "D48" "0.4683"
Set-TextValue "E48" "  +0.37%  "
Set-TextValue "E49" "  +0.01%  "
Set-TextValue "D50" "1.608"
Set-TextValue "E50" "  -1.95%  "
Set-TextValue "D51" "63.68"
Set-TextValue "E51" "  -1.53%  "
$excel.CutCopyMode = 0

